# Common: Build edit works
# Adds the "lab.build.*" translation rows to the Import sheet of the
# translations workbook (rows 510-519), mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

function Set-ImportRow($r, $k, $t) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = "cs"
    $cellA.WrapText = $true
    $cellA.Font.Size = 10

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $k
    $cellB.WrapText = $true
    $cellB.Font.Size = 10

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = $t
    $cellC.WrapText = $true
    $cellC.Font.Size = 10
}

Set-ImportRow 510 "lab.build.index.title" "Detail buildu"
Set-ImportRow 511 "lab.build.index.preview.title" "Detail buildu"
Set-ImportRow 512 "lab.build.index.preview.subtitle" "Zde můžete spravovat vybraný build."
Set-ImportRow 513 "lab.build.button.edit" "Upravit"
Set-ImportRow 514 "lab.build.preview.name" "Jméno"
Set-ImportRow 515 "lab.build.edit.title" "Editace buildu"
Set-ImportRow 516 "lab.build.edit.subtitle" "Místo, kde je možné upravit build."
Set-ImportRow 517 "lab.build.update.submit" "Aktualizovat"
Set-ImportRow 518 "lab.build.update.success" "Build [{{data.name}}] byl aktualizován."
Set-ImportRow 519 "lab.build.link.button" "Zpět"

# Match the author's final scroll position / selection on the sheet.
$ws.Activate()
$ws.Range("B514").Select() | Out-Null
